$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 3383.1333
$ws.Range("I11").Value = 3383.1333
$ws.Range("K11").Value = 3383.1333
$ws.Range("M11").Value = -3243.1333
$ws.Range("H33").Value = 276.07144
$ws.Range("I33").Value = 297.72726
$ws.Range("J33").Value = 196.66667
$ws.Range("K33").Value = 297.72726
$ws.Range("L33").Value = 196.66667
$ws.Range("M33").Value = -68.72726
$ws.Range("N33").Value = -654.6666700000001
$ws.Range("H38").Value = 130.41667
$ws.Range("I38").Value = 130.41667
$ws.Range("K38").Value = 391.25001
$ws.Range("M38").Value = -19.25001000000003
$ws.Range("H46").Value = 4258.5
$ws.Range("I46").Value = 517
$ws.Range("K46").Value = 1551
$ws.Range("M46").Value = -1432
$ws.Range("H60").Value = 4258.5
$ws.Range("I60").Value = 517
$ws.Range("K60").Value = 1551
$ws.Range("M60").Value = -1067
$ws.Range("H62").Value = 211115020
$ws.Range("I62").Value = 211115020
$ws.Range("K62").Value = 211115020
$ws.Range("M62").Value = -211114396
$ws.Range("H65").Value = 211115020
$ws.Range("I65").Value = 211115020
$ws.Range("K65").Value = 1055575100
$ws.Range("M65").Value = -1055571980
$ws.Range("H98").Value = 1138
$ws.Range("I98").Value = 982.9167
$ws.Range("J98").Value = 2999
$ws.Range("K98").Value = 982.9167
$ws.Range("L98").Value = 2999
$ws.Range("M98").Value = 515.0833
$ws.Range("N98").Value = -5995
$ws.Range("H106").Value = 3358.4
$ws.Range("I106").Value = 3265.6667
$ws.Range("K106").Value = 3265.6667
$ws.Range("M106").Value = -2634.6667
$ws.Range("H122").Value = 1138
$ws.Range("I122").Value = 982.9167
$ws.Range("J122").Value = 2999
$ws.Range("K122").Value = 2948.7501
$ws.Range("L122").Value = 8997
$ws.Range("M122").Value = -498.7501000000002
$ws.Range("N122").Value = -13897

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4843535.5
$ws.Range("I2").Value = 6134852
$ws.Range("J2").Value = 1099
$ws.Range("K2").Value = 6134852
$ws.Range("L2").Value = 1099
$ws.Range("M2").Value = -6134739
$ws.Range("N2").Value = -1325
$ws.Range("H63").Value = 7813
$ws.Range("J63").Value = 9920.5
$ws.Range("L63").Value = 9920.5
$ws.Range("N63").Value = -11292.5
$ws.Range("H66").Value = 7813
$ws.Range("J66").Value = 9920.5
$ws.Range("L66").Value = 49602.5
$ws.Range("N66").Value = -56466.5
$ws.Range("H97").Value = 2058184.2
$ws.Range("I97").Value = 2179245.2
$ws.Range("J97").Value = 145
$ws.Range("K97").Value = 2179245.2
$ws.Range("L97").Value = 145
$ws.Range("M97").Value = -2178749.2
$ws.Range("N97").Value = -1137
$ws.Range("H102").Value = 14766328
$ws.Range("I102").Value = 1897.6428
$ws.Range("J102").Value = 83667000
$ws.Range("K102").Value = 1897.6428
$ws.Range("L102").Value = 83667000
$ws.Range("M102").Value = -275.6428000000001
$ws.Range("N102").Value = -83670244
$ws.Range("H116").Value = 4843535.5
$ws.Range("I116").Value = 6134852
$ws.Range("J116").Value = 1099
$ws.Range("K116").Value = 6134852
$ws.Range("L116").Value = 1099
$ws.Range("M116").Value = -6132558
$ws.Range("N116").Value = -5687
$ws.Range("H122").Value = 1792.5
$ws.Range("I122").Value = 750
$ws.Range("J122").Value = 2001
$ws.Range("K122").Value = 2250
$ws.Range("L122").Value = 6003
$ws.Range("M122").Value = 200
$ws.Range("N122").Value = -10903

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4843535.5
$ws.Range("I3").Value = 6134852
$ws.Range("J3").Value = 1099
$ws.Range("K3").Value = 6134852
$ws.Range("L3").Value = 1099
$ws.Range("M3").Value = -6134738
$ws.Range("N3").Value = -1327
$ws.Range("H80").Value = 52386.5
$ws.Range("J80").Value = 80325.92
$ws.Range("L80").Value = 80325.92
$ws.Range("N80").Value = -82321.92
$ws.Range("H83").Value = 52386.5
$ws.Range("J83").Value = 80325.92
$ws.Range("L83").Value = 401629.6
$ws.Range("N83").Value = -411613.6
$ws.Range("H105").Value = 100027780
$ws.Range("I105").Value = 125033600
$ws.Range("J105").Value = 4505
$ws.Range("K105").Value = 125033600
$ws.Range("L105").Value = 4505
$ws.Range("M105").Value = -125031853
$ws.Range("N105").Value = -7999
$ws.Range("H107").Value = 1771.5
$ws.Range("I107").Value = 1298.6666
$ws.Range("J107").Value = 2007.9166
$ws.Range("K107").Value = 1298.6666
$ws.Range("L107").Value = 2007.9166
$ws.Range("M107").Value = 621.3334
$ws.Range("N107").Value = -5847.9166
$ws.Range("H134").Value = 4732.457
$ws.Range("I134").Value = 3207.1853
$ws.Range("K134").Value = 9621.555899999999
$ws.Range("M134").Value = -7086.555899999999

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 17396.285
$ws.Range("I41").Value = 3569.6667
$ws.Range("J41").Value = 27766.25
$ws.Range("K41").Value = 3569.6667
$ws.Range("L41").Value = 27766.25
$ws.Range("M41").Value = -3141.6667
$ws.Range("N41").Value = -28622.25
$ws.Range("H105").Value = 2759.5715
$ws.Range("I105").Value = 1552.8334
$ws.Range("J105").Value = 10000
$ws.Range("K105").Value = 1552.8334
$ws.Range("L105").Value = 10000
$ws.Range("M105").Value = 194.1666
$ws.Range("N105").Value = -13494
$ws.Range("H132").Value = 70946.89
$ws.Range("I132").Value = 9117.799999999999
$ws.Range("K132").Value = 27353.4
$ws.Range("M132").Value = -24823.4

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 739.25
$ws.Range("I5").Value = 677.3333
$ws.Range("J5").Value = 925
$ws.Range("K5").Value = 2031.9999
$ws.Range("L5").Value = 2775
$ws.Range("M5").Value = -1919.9999
$ws.Range("N5").Value = -2999
$ws.Range("H34").Value = 1723.6
$ws.Range("J34").Value = 6340.2
$ws.Range("L34").Value = 19020.6
$ws.Range("N34").Value = -19188.6
$ws.Range("H101").Value = 3163
$ws.Range("J101").Value = 4000
$ws.Range("L101").Value = 12000
$ws.Range("N101").Value = -16868
$ws.Range("H102").Value = 4666.6665
$ws.Range("I102").Value = 4666.6665
$ws.Range("K102").Value = 13999.9995
$ws.Range("M102").Value = -11565.9995
$ws.Range("H108").Value = 1615.4
$ws.Range("I108").Value = 781.75
$ws.Range("J108").Value = 4950
$ws.Range("K108").Value = 2345.25
$ws.Range("L108").Value = 14850
$ws.Range("M108").Value = 534.75
$ws.Range("N108").Value = -20610
$ws.Range("H110").Value = 4263.5
$ws.Range("I110").Value = 1027
$ws.Range("K110").Value = 3081
$ws.Range("M110").Value = 1009
$ws.Range("H135").Value = 739.25
$ws.Range("I135").Value = 677.3333
$ws.Range("J135").Value = 925
$ws.Range("K135").Value = 6095.9997
$ws.Range("L135").Value = 8325
$ws.Range("M135").Value = -3560.9997
$ws.Range("N135").Value = -13395

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4950.75
$ws.Range("I132").Value = 3490.6333
$ws.Range("J132").Value = 8079.5713
$ws.Range("K132").Value = 10471.8999
$ws.Range("L132").Value = 24238.7139
$ws.Range("M132").Value = -7941.8999
$ws.Range("N132").Value = -29298.7139

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 9350
$ws.Range("J4").Value = 8700
$ws.Range("L4").Value = 8700
$ws.Range("N4").Value = -8926
$ws.Range("H28").Value = 9350
$ws.Range("J28").Value = 8700
$ws.Range("L28").Value = 8700
$ws.Range("N28").Value = -9164
$ws.Range("H34").Value = 12000
$ws.Range("J34").Value = 12000
$ws.Range("L34").Value = 12000
$ws.Range("N34").Value = -12344
$ws.Range("H35").Value = 3459.7144
$ws.Range("I35").Value = 1452.8334
$ws.Range("J35").Value = 15501
$ws.Range("K35").Value = 1452.8334
$ws.Range("L35").Value = 15501
$ws.Range("M35").Value = -1116.8334
$ws.Range("N35").Value = -16173
$ws.Range("H37").Value = 9350
$ws.Range("J37").Value = 8700
$ws.Range("L37").Value = 8700
$ws.Range("N37").Value = -8914
$ws.Range("H48").Value = 60000
$ws.Range("I48").Value = 60000
$ws.Range("K48").Value = 60000
$ws.Range("M48").Value = -59339
$ws.Range("H100").Value = 13891957
$ws.Range("I100").Value = 31252436
$ws.Range("K100").Value = 31252436
$ws.Range("M100").Value = -31251895
$ws.Range("H136").Value = 11120243
$ws.Range("I136").Value = 16673940
$ws.Range("K136").Value = 50021820
$ws.Range("M136").Value = -50019270

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H68").Value = 36931.668
$ws.Range("J68").Value = 36931.668
$ws.Range("L68").Value = 36931.668
$ws.Range("N68").Value = -38553.668
$ws.Range("H71").Value = 36931.668
$ws.Range("J71").Value = 36931.668
$ws.Range("L71").Value = 110795.004
$ws.Range("N71").Value = -118907.004
$ws.Range("H81").Value = 13825.193
$ws.Range("I81").Value = 6647.1
$ws.Range("K81").Value = 13294.2
$ws.Range("M81").Value = -12233.2
$ws.Range("H84").Value = 13825.193
$ws.Range("I84").Value = 6647.1
$ws.Range("K84").Value = 66471
$ws.Range("M84").Value = -61167
$ws.Range("H107").Value = 1933.8572
$ws.Range("I107").Value = 1394.3636
$ws.Range("J107").Value = 2527.3
$ws.Range("K107").Value = 4183.0908
$ws.Range("L107").Value = 7581.900000000001
$ws.Range("M107").Value = -2263.0908
$ws.Range("N107").Value = -11421.9
$ws.Range("H132").Value = 10931
$ws.Range("I132").Value = 9399.666999999999
$ws.Range("K132").Value = 28199.001
$ws.Range("M132").Value = -25669.001
